$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update accuracy values changed by the re-run of Indonesian PoS and PUD tests
$ws.Range("C2").Value = 0.9787187739463602
$ws.Range("D2").Value = 0.810360153256705
$ws.Range("E2").Value = 0.7337164750957854
$ws.Range("F2").Value = 0.8548934865900383
$ws.Range("G2").Value = 0.8065409961685823
$ws.Range("H2").Value = 0.6467310344827586
$ws.Range("I2").Value = 0.557391570881226
$ws.Range("L2").Value = 0.7021180076628353
$ws.Range("M2").Value = 0.6832613026819924
$ws.Range("N2").Value = 0.6263785440613027
$ws.Range("O2").Value = 0.5810298850574712
$ws.Range("P2").Value = 0.5382283524904214
$ws.Range("Q2").Value = 0.6748842911877395
$ws.Range("R2").Value = 0.6471846743295019
$ws.Range("S2").Value = 0.7383662835249042
$ws.Range("T2").Value = 0.6823540229885058
$ws.Range("U2").Value = 0.5052137931034483
$ws.Range("C3").Value = 0.8256851196000076
$ws.Range("D3").Value = 0.9869131266453287
$ws.Range("E3").Value = 0.8633359216681502
$ws.Range("F3").Value = 0.8864983617734513
$ws.Range("G3").Value = 0.7844169617999659
$ws.Range("H3").Value = 0.5815041381792011
$ws.Range("I3").Value = 0.5821480653775497
$ws.Range("L3").Value = 0.8327115016761046
$ws.Range("M3").Value = 0.626048749076722
$ws.Range("N3").Value = 0.5896668623700309
$ws.Range("O3").Value = 0.5204068104770743
$ws.Range("P3").Value = 0.4949527471070624
$ws.Range("Q3").Value = 0.6242684797636409
$ws.Range("R3").Value = 0.6516353856934528
$ws.Range("S3").Value = 0.7544175299710233
$ws.Range("T3").Value = 0.7061608681653757
$ws.Range("U3").Value = 0.402435559932577
$ws.Range("L4").Value = 0.7619742093951489
$ws.Range("L5").Value = 0.7897283588066475
$ws.Range("L6").Value = 0.6918146859549533
$ws.Range("C7").Value = 0.5348817848817848
$ws.Range("D7").Value = 0.4408924408924409
$ws.Range("E7").Value = 0.3677156177156177
$ws.Range("F7").Value = 0.6205461205461206
$ws.Range("G7").Value = 0.6103896103896104
$ws.Range("H7").Value = 0.9638694638694638
$ws.Range("I7").Value = 0.5282217782217782
$ws.Range("L7").Value = 0.571012321012321
$ws.Range("M7").Value = 0.5994838494838495
$ws.Range("N7").Value = 0.5815850815850816
$ws.Range("O7").Value = 0.6424408924408924
$ws.Range("P7").Value = 0.6506826506826506
$ws.Range("Q7").Value = 0.5231435231435232
$ws.Range("R7").Value = 0.5092407592407593
$ws.Range("S7").Value = 0.6644189144189144
$ws.Range("T7").Value = 0.5777555777555777
$ws.Range("U7").Value = 0.4756909756909757
$ws.Range("L8").Value = 0.6311166875784191
$ws.Range("L9").Value = 0.6224800645103485
$ws.Range("L10").Value = 0.4883603966087082
$ws.Range("C11").Value = 0.6864176570458405
$ws.Range("D11").Value = 0.8378607809847198
$ws.Range("E11").Value = 0.7747877758913413
$ws.Range("F11").Value = 0.8368421052631579
$ws.Range("G11").Value = 0.6973684210526315
$ws.Range("H11").Value = 0.6874363327674023
$ws.Range("I11").Value = 0.5953310696095077
$ws.Range("L11").Value = 0.9286078098471986
$ws.Range("M11").Value = 0.749660441426146
$ws.Range("N11").Value = 0.6971986417657046
$ws.Range("O11").Value = 0.6166383701188455
$ws.Range("P11").Value = 0.582258064516129
$ws.Range("Q11").Value = 0.6729202037351443
$ws.Range("R11").Value = 0.6131578947368421
$ws.Range("S11").Value = 0.7716468590831919
$ws.Range("T11").Value = 0.6694397283531409
$ws.Range("U11").Value = 0.5459252971137522
$ws.Range("C12").Value = 0.77695730379627
$ws.Range("D12").Value = 0.851415317618101
$ws.Range("E12").Value = 0.7652182145223895
$ws.Range("F12").Value = 0.8307299062766259
$ws.Range("G12").Value = 0.8305879011644419
$ws.Range("H12").Value = 0.7406986651519455
$ws.Range("I12").Value = 0.5973681719208558
$ws.Range("L12").Value = 0.7626147874656821
$ws.Range("M12").Value = 0.9675754993846445
$ws.Range("N12").Value = 0.8019975385780554
$ws.Range("O12").Value = 0.6774590551926536
$ws.Range("P12").Value = 0.7219066553062576
$ws.Range("Q12").Value = 0.7824008330966582
$ws.Range("R12").Value = 0.6713055003313453
$ws.Range("S12").Value = 0.8437943765975575
$ws.Range("T12").Value = 0.7154690902205812
$ws.Range("U12").Value = 0.5933920287797028
$ws.Range("L13").Value = 0.6957003364240584
$ws.Range("L14").Value = 0.5687090178382571
$ws.Range("C15").Value = 0.3985592357685381
$ws.Range("D15").Value = 0.3056142823584684
$ws.Range("E15").Value = 0.2662281732049174
$ws.Range("F15").Value = 0.4128885756792733
$ws.Range("G15").Value = 0.4604964372406233
$ws.Range("H15").Value = 0.5389554459321901
$ws.Range("I15").Value = 0.3994205622112599
$ws.Range("L15").Value = 0.3534570511314697
$ws.Range("M15").Value = 0.4625322997416021
$ws.Range("N15").Value = 0.5369978858350951
$ws.Range("O15").Value = 0.58656330749354
$ws.Range("P15").Value = 0.9791715605669093
$ws.Range("Q15").Value = 0.5359799545846058
$ws.Range("R15").Value = 0.3555712160363323
$ws.Range("S15").Value = 0.4778795709028267
$ws.Range("T15").Value = 0.487902278599953
$ws.Range("U15").Value = 0.3990290501918409
$ws.Range("L16").Value = 0.6893895065340355
$ws.Range("L17").Value = 0.6246587956720492
$ws.Range("L18").Value = 0.6235713813833652
$ws.Range("L19").Value = 0.288418206399279
$ws.Range("L20").Value = 0.2045516120292604
$ws.Range("C22").Value = 0.5911523355305596
$ws.Range("D22").Value = 0.6335096123274141
$ws.Range("E22").Value = 0.5881901813559488
$ws.Range("F22").Value = 0.6134407436271602
$ws.Range("G22").Value = 0.6391701995579466
$ws.Range("H22").Value = 0.5889554720406242
$ws.Range("I22").Value = 0.5219256242580397
$ws.Range("L22").Value = 0.6056881626712747
$ws.Range("M22").Value = 0.5945463845998242
$ws.Range("N22").Value = 0.5886170237731682
$ws.Range("O22").Value = 0.5353549072451216
$ws.Range("P22").Value = 0.5431152791252254
$ws.Range("Q22").Value = 0.5760540638225131
$ws.Range("R22").Value = 0.5552780983803525
$ws.Range("S22").Value = 0.6596798515668796
$ws.Range("T22").Value = 0.6194941661509865
$ws.Range("U22").Value = 0.480827564021314
$ws.Range("C27").Value = 0.8355683181540132
$ws.Range("D27").Value = 0.843948512539462
$ws.Range("E27").Value = 0.7900824301343616
$ws.Range("F27").Value = 0.8442304659166722
$ws.Range("G27").Value = 0.8036671029611815
$ws.Range("H27").Value = 0.6350996109242726
$ws.Range("I27").Value = 0.5649516931878795
$ws.Range("L27").Value = 0.7556693526991379
$ws.Range("M27").Value = 0.7031024411549016
$ws.Range("N27").Value = 0.6422269487762176
$ws.Range("O27").Value = 0.5694710925400618
$ws.Range("P27").Value = 0.5369960812746044
$ws.Range("Q27").Value = 0.6686089213326457
$ws.Range("R27").Value = 0.6857143806573156
$ws.Range("S27").Value = 0.7533063541494186
$ws.Range("T27").Value = 0.6938093982886983
$ws.Range("U27").Value = 0.4990957398990363
$ws.Range("C28").Value = 0.5481104333791165
$ws.Range("D28").Value = 0.5768306641626223
$ws.Range("E28").Value = 0.5515820772042
$ws.Range("F28").Value = 0.620568718708227
$ws.Range("G28").Value = 0.6150689583260652
$ws.Range("H28").Value = 0.6871789478225553
$ws.Range("I28").Value = 0.5429234779317477
$ws.Range("L28").Value = 0.5782423674274493
$ws.Range("M28").Value = 0.6164785463784395
$ws.Range("N28").Value = 0.599478359641083
$ws.Range("O28").Value = 0.5758638019390909
$ws.Range("P28").Value = 0.5828923246897985
$ws.Range("Q28").Value = 0.5460820483632336
$ws.Range("R28").Value = 0.5555550215399206
$ws.Range("S28").Value = 0.6759196328225643
$ws.Range("T28").Value = 0.5933781872420358
$ws.Range("U28").Value = 0.4688687235335968
$ws.Range("C29").Value = 0.6061951661742485
$ws.Range("D29").Value = 0.6290105436377385
$ws.Range("E29").Value = 0.579735806092556
$ws.Range("F29").Value = 0.6295883380069857
$ws.Range("G29").Value = 0.6466451715177792
$ws.Range("H29").Value = 0.6422907871590777
$ws.Range("I29").Value = 0.5460161349399134
$ws.Range("L29").Value = 0.6139741398787006
$ws.Range("M29").Value = 0.6246909267948835
$ws.Range("N29").Value = 0.6666421500649489
$ws.Range("O29").Value = 0.6279718323457103
$ws.Range("P29").Value = 0.6827780440518729
$ws.Range("Q29").Value = 0.6496462751434467
$ws.Range("R29").Value = 0.5349718102811883
$ws.Range("S29").Value = 0.6694420657237685
$ws.Range("T29").Value = 0.6086274008577895
$ws.Range("U29").Value = 0.5008237058357512
$ws.Range("L30").Value = 0.4352999988709884
$ws.Range("C32").Value = 0.5929022774997944
$ws.Range("D32").Value = 0.6373332383808185
$ws.Range("E32").Value = 0.5910065863193952
$ws.Range("F32").Value = 0.6119860205358545
$ws.Range("G32").Value = 0.6402093413874497
$ws.Range("H32").Value = 0.5827380061654928
$ws.Range("I32").Value = 0.5177308380323076
$ws.Range("L32").Value = 0.5957964647190691
$ws.Range("M32").Value = 0.5863908709539685
$ws.Range("N32").Value = 0.584587569968733
$ws.Range("O32").Value = 0.5306908397458147
$ws.Range("P32").Value = 0.5410116636431034
$ws.Range("Q32").Value = 0.5721426361943347
$ws.Range("R32").Value = 0.5368769454907912
$ws.Range("S32").Value = 0.6397417685676944
$ws.Range("T32").Value = 0.6132720848283925
$ws.Range("U32").Value = 0.4764433178205735
$ws.Range("C37").Value = 0.8234993659411382
$ws.Range("D37").Value = 0.6519068856037634
$ws.Range("E37").Value = 0.6240810970156861
$ws.Range("F37").Value = 0.6579814682486173
$ws.Range("C38").Value = 0.5824321703560462
$ws.Range("D38").Value = 0.6027815977272507
$ws.Range("E38").Value = 0.5841590162023291
$ws.Range("F38").Value = 0.5734303912845294
$ws.Range("C39").Value = 0.6182350050858616
$ws.Range("D39").Value = 0.600760353992564
$ws.Range("E39").Value = 0.6503458456801725
$ws.Range("F39").Value = 0.5784662456746243
$ws.Range("D40").Value = 0.4062382412322479
$ws.Range("C42").Value = 0.5450835351191705
$ws.Range("D42").Value = 0.552968493609525
$ws.Range("E42").Value = 0.533837672908197
$ws.Range("F42").Value = 0.6032927017359236

# The bold/underlined "row max" highlight moved from D27 to F27
$ws.Range("D27").Copy()
$ws.Range("F27").PasteSpecial(-4122)  # xlPasteFormats (from old D27, which was the max-highlight style)
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)  # xlPasteFormats (reset D27 back to the plain style)
$excel.CutCopyMode = 0
